# Add new JVerb (JUCE reverb) effect
# - Rename the existing "Reverb" effect (rows 79-86) to "Mverb"
# - Add a brand-new "JuceVerb" effect block in rows 88-93

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename "Reverb" -> "Mverb" for the existing effect block (rows 79-86) ---
$ws.Range("A79:A86").Value = "Mverb"

# --- New "JuceVerb" effect rows (row 87 intentionally left blank, matching the gap before row 88) ---
# Fill parameter name/variable/label columns (B, C, D, K) first, then the G/H/I numeric
# bounds, and only then the effect-name column (A) last for every new row.

# Row 88: Room Size
$ws.Range("B88").Value = "Room Size"
$ws.Range("C88").Value = "roomSize"
$ws.Range("D88").Value = "Room Size"
$ws.Range("K88").Value = "roomSize"

# Row 89: Damping
$ws.Range("B89").Value = "Damping"
$ws.Range("C89").Value = "damping"
$ws.Range("D89").Value = "Damping"
$ws.Range("K89").Value = "damping"

# Row 90: Wet Level
$ws.Range("B90").Value = "Wet Level"
$ws.Range("C90").Value = "wetLevel"
$ws.Range("D90").Value = "Wet Level"
$ws.Range("K90").Value = "wetLevel"

# Row 91: Dry Level
$ws.Range("B91").Value = "Dry Level"
$ws.Range("C91").Value = "dryLevel"
$ws.Range("D91").Value = "Dry Level"
$ws.Range("K91").Value = "dryLevel"

# Row 92: Width
$ws.Range("B92").Value = "Width"
$ws.Range("C92").Value = "width"
$ws.Range("D92").Value = "Width"
$ws.Range("K92").Value = "width"

# Row 93: Freeze Mode (enum-style row: E = enum count, I = default index)
$ws.Range("B93").Value = "Freeze Mode"
$ws.Range("C93").Value = "freezeMode"
$ws.Range("D93").Value = "Freeze"
$ws.Range("K93").Value = "freezeMode"

# Numeric bounds / defaults
$ws.Range("G88").Value = 0
$ws.Range("H88").Value = 1
$ws.Range("I88").Value = 0.5

$ws.Range("G89").Value = 0
$ws.Range("H89").Value = 1
$ws.Range("I89").Value = 0.5

$ws.Range("G90").Value = 0
$ws.Range("H90").Value = 1
$ws.Range("I90").Value = 0.33

$ws.Range("G91").Value = 0
$ws.Range("H91").Value = 1
$ws.Range("I91").Value = 0.4

$ws.Range("G92").Value = 0
$ws.Range("H92").Value = 1
$ws.Range("I92").Value = 1

$ws.Range("E93").Value = 2
$ws.Range("I93").Value = 0

# Effect-name column filled in last across the whole new block
$ws.Range("A88:A93").Value = "JuceVerb"

# Match the saved selection state from the edit: A89:A93 selected, A89 active
$ws.Range("A89:A93").Select()
